# Add a new test-case row (uid=11, nric="T0193322F", password="password")
# to the User worksheet, right after the existing last row (row 12).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "T0193322F"
$ws.Range("C13").Value = "password"
